# Apply TestData updates for the AI-driven test.
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("TestData")
$wsLists = $wb.Worksheets.Item("Lists")

# Update URL values from "new/..." to "#/..."
$wsData.Range("B3").Value = "#/login"
$wsData.Range("B4").Value = "#/catalog"

# Row 6: B6 becomes "security", C6 becomes "NewUI"
$wsData.Range("B6").Value = "security"
$wsData.Range("C6").Value = "NewUI"

# Update the active selection from B7 to B6
$wsData.Activate()
$wsData.Range("B6").Select()
